$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new rows of order data (context data) beneath the existing rows.
$ws.Range("A30").Value = "183090-0"
$ws.Range("B30").Value = "Clio - Greek Yogurt Bar Vanilla"
$ws.Range("C30").Value = "'1"
$ws.Range("D30").Value = "'15.45"
$ws.Range("E30").Value = "'15.45"

$ws.Range("A31").Value = "183096-7"
$ws.Range("B31").Value = "Clio - Greek Yogurt Bar Strawberry"
$ws.Range("C31").Value = "'1"
$ws.Range("D31").Value = "'15.45"
$ws.Range("E31").Value = "'15.45"

# Re-apply the plain default formatting used throughout the rest of the sheet
# so the quote-prefix ("stored as text") style introduced above by the
# leading-apostrophe text entry doesn't leave these new cells on a different
# style than the rest of the table.
$ws.Range("A1:E1").Copy()
$ws.Range("A30:E31").PasteSpecial(-4122)
$excel.CutCopyMode = $false
